$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Only the IAM module (row 2) keeps Runmode = "Y".
# All other modules (rows 3-7) should now run with Runmode = "N".
$ws.Range("C3:C7").Value = "N"
